$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I13").Value = "sv"
$ws.Range("J13").Value = "Statement-opinion"
$ws.Range("I18").Value = "sd"
$ws.Range("J18").Value = "Statement-non-opinion"
$ws.Range("I21").Value = "sd"
$ws.Range("J21").Value = "Statement-non-opinion"
$ws.Range("I22").Value = "sd"
$ws.Range("J22").Value = "Statement-non-opinion"
$ws.Range("I43").Value = "sd"
$ws.Range("J43").Value = "Statement-non-opinion"
$ws.Range("I46").Value = "sv"
$ws.Range("J46").Value = "Statement-opinion"
$ws.Range("I49").Value = "sv"
$ws.Range("J49").Value = "Statement-opinion"
$ws.Range("I62").Value = "ba"
$ws.Range("J62").Value = "Appreciation"
$ws.Range("I71").Value = "sd"
$ws.Range("J71").Value = "Statement-non-opinion"
$ws.Range("I98").Value = "sv"
$ws.Range("J98").Value = "Statement-opinion"
$ws.Range("I124").Value = "%"
$ws.Range("J124").Value = "Uninterpretable"
$ws.Range("I128").Value = "sd"
$ws.Range("J128").Value = "Statement-non-opinion"
$ws.Range("I129").Value = "sd"
$ws.Range("J129").Value = "Statement-non-opinion"
$ws.Range("I134").Value = "ba"
$ws.Range("J134").Value = "Appreciation"
$ws.Range("I135").Value = "sd"
$ws.Range("J135").Value = "Statement-non-opinion"
$ws.Range("I136").Value = "sd"
$ws.Range("J136").Value = "Statement-non-opinion"
$ws.Range("I141").Value = "sv"
$ws.Range("J141").Value = "Statement-opinion"
$ws.Range("I148").Value = "%"
$ws.Range("J148").Value = "Uninterpretable"
$ws.Range("I151").Value = "sd"
$ws.Range("J151").Value = "Statement-non-opinion"
$ws.Range("I153").Value = "sd"
$ws.Range("J153").Value = "Statement-non-opinion"
$ws.Range("I170").Value = "aa"
$ws.Range("J170").Value = "Agree/Accept"
$ws.Range("I171").Value = "aa"
$ws.Range("J171").Value = "Agree/Accept"
$ws.Range("I172").Value = "sd"
$ws.Range("J172").Value = "Statement-non-opinion"
$ws.Range("I174").Value = "b"
$ws.Range("J174").Value = "Acknowledge (Backchannel)"
$ws.Range("I184").Value = "aa"
$ws.Range("J184").Value = "Agree/Accept"
$ws.Range("I194").Value = "sv"
$ws.Range("J194").Value = "Statement-opinion"
$ws.Range("I203").Value = "ba"
$ws.Range("J203").Value = "Appreciation"
$ws.Range("I205").Value = "sd"
$ws.Range("J205").Value = "Statement-non-opinion"
$ws.Range("I216").Value = "aa"
$ws.Range("J216").Value = "Agree/Accept"
$ws.Range("I229").Value = "sv"
$ws.Range("J229").Value = "Statement-opinion"
$ws.Range("I234").Value = "qy"
$ws.Range("J234").Value = "Yes-No-Question"
$ws.Range("I237").Value = "%"
$ws.Range("J237").Value = "Uninterpretable"
$ws.Range("I241").Value = "sd"
$ws.Range("J241").Value = "Statement-non-opinion"
